# TA-01/ Admin Scenario modified with verification of deleted user
#
# Administration sheet ("Administration"):
#  - TS22 (row 43): the FILTER_TABLE element is no longer found, so the
#    step result flips from PASS to FAIL and the logged message changes.
#  - New steps TS29-TS35 are appended (rows 55-61): open the user
#    dropdown, log out, and log back in with the just-deleted
#    AutomationUser credentials.
#  - TS24's "login error" check (already used earlier, on the Settings
#    sheet) is reused as the final new row (62) to verify the deleted
#    user can no longer log in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Administration")

# --- 1. FILTER_TABLE step now fails ---------------------------------
$ws.Range("F43").Value = "FAIL"
$ws.Range("G43").Value = "Element FILTER_TABLE not found"

# --- 2. Append the new verification steps (rows 55-62) --------------
# Copy the formatting (style) of the last existing data row onto the
# new rows first, so the A:E columns keep the same text-formatted
# style (s="1") the rest of the table uses.
$ws.Range("A54:E54").Copy()
$ws.Range("A55:E62").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rows = @(
  @{ Row=55; A="TS29"; B="xpath";    C="USER_DROPDOWN";    D="click";          E="";                                                                       F="PASS"; G="clicked on USER_DROPDOWN" },
  @{ Row=56; A="TS30"; B="";         C="";                 D="wait";           E="2000";                                                                   F="PASS"; G="Waiting for 2000 millisecond  text box to be present" },
  @{ Row=57; A="TS31"; B="linkText"; C="LOGOUT_LINK";      D="click";          E="";                                                                       F="PASS"; G="clicked on LOGOUT_LINK" },
  @{ Row=58; A="TS32"; B="";         C="";                 D="verifyTitle";    E="oneView";                                                                F="PASS"; G="TITLE matches : oneView" },
  @{ Row=59; A="TS33"; B="id";       C="USERNAME_TEXTBOX"; D="sendKeys";       E="AutomationUser";                                                          F="PASS"; G="Typed AutomationUser into USERNAME_TEXTBOX text box" },
  @{ Row=60; A="TS34"; B="id";       C="PASSWORD_TEXTBOX"; D="sendKeys";       E="Auto123";                                                                 F="PASS"; G="Typed Auto123 into PASSWORD_TEXTBOX text box" },
  @{ Row=61; A="TS35"; B="id";       C="SUBMIT_BUTTON";    D="click";          E="";                                                                       F="PASS"; G="clicked on SUBMIT_BUTTON" },
  @{ Row=62; A="TS24"; B="xpath";    C="LOGIN_ERRORMSG";   D="getTextContent"; E="Login Error: The username/password you provided were incorrect.";        F="PASS"; G="Element present : Login Error: The username/password you provided were incorrect." }
)

foreach ($r in $rows) {
  $n = $r.Row
  $ws.Range("A$n").Value = $r.A
  if ($r.B -ne "") { $ws.Range("B$n").Value = $r.B }
  if ($r.C -ne "") { $ws.Range("C$n").Value = $r.C }
  $ws.Range("D$n").Value = $r.D
  if ($r.E -ne "") { $ws.Range("E$n").Value = $r.E }
  $ws.Range("F$n").Value = $r.F
  $ws.Range("G$n").Value = $r.G
}

# Row 61's (E) cell keeps the distinct "Hyperlink look" formatting used
# sporadically on other blank Data cells elsewhere in this column.
$ws.Range("E5").Copy()
$ws.Range("E61").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Restore the view/selection state on the affected sheets -----
$ws.Activate()
$ws.Range("I52").Select()

$ws2 = $wb.Worksheets.Item("Settings")
$ws2.Activate()
$ws2.Range("A25:E25").Select()

$ws.Activate()
